{"js": "// Edit: 1) split the title paragraph's run so the \"CITTA' DI IMPERIA\" picture\n//          is its own run that precedes a separate run holding the text, and\n//       2) fix \"informarLe\" -> \"informarLa\" in the body paragraph.\n//\n// `context` (alias `ctx`) is the Word.RequestContext; this is the body of\n// `async (context) => { ... }`.\n\n// --- 1) Re-shape the first paragraph (picture + \"CITTA' DI IMPERIA\") -------\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst titlePara = paragraphs.items[0];\n\n// The picture keeps its original anchor/position/size \u2014 only the run\n// structure changes: the drawing becomes its own run, immediately followed\n// by a run that just carries the \"CITTA' DI IMPERIA\" text.\nconst drawingXml =\n  '<w:drawing>' +\n    '<wp:anchor behindDoc=\"0\" distT=\"0\" distB=\"0\" distL=\"0\" distR=\"0\" simplePos=\"0\" locked=\"0\" layoutInCell=\"1\" allowOverlap=\"1\" relativeHeight=\"2\">' +\n      '<wp:simplePos x=\"0\" y=\"0\"/>' +\n      '<wp:positionH relativeFrom=\"column\"><wp:posOffset>786765</wp:posOffset></wp:positionH>' +\n      '<wp:positionV relativeFrom=\"paragraph\"><wp:posOffset>36195</wp:posOffset></wp:positionV>' +\n      '<wp:extent cx=\"704215\" cy=\"1000125\"/>' +\n      '<wp:effectExtent l=\"0\" t=\"0\" r=\"0\" b=\"0\"/>' +\n      '<wp:wrapTopAndBottom/>' +\n      '<wp:docPr id=\"1\" name=\"Picture\" descr=\"\"/>' +\n      '<wp:cNvGraphicFramePr><a:graphicFrameLocks xmlns:a=\"http://schemas.openxmlformats.org/drawingml/2006/main\" noChangeAspect=\"1\"/></wp:cNvGraphicFramePr>' +\n      '<a:graphic xmlns:a=\"http://schemas.openxmlformats.org/drawingml/2006/main\">' +\n        '<a:graphicData uri=\"http://schemas.openxmlformats.org/drawingml/2006/picture\">' +\n          '<pic:pic xmlns:pic=\"http://schemas.openxmlformats.org/drawingml/2006/picture\">' +\n            '<pic:nvPicPr>' +\n              '<pic:cNvPr id=\"1\" name=\"Picture\" descr=\"\"/>' +\n              '<pic:cNvPicPr><a:picLocks noChangeAspect=\"1\" noChangeArrowheads=\"1\"/></pic:cNvPicPr>' +\n            '</pic:nvPicPr>' +\n            '<pic:blipFill><a:blip r:embed=\"rId2\"/><a:stretch><a:fillRect/></a:stretch></pic:blipFill>' +\n            '<pic:spPr bwMode=\"auto\">' +\n              '<a:xfrm><a:off x=\"0\" y=\"0\"/><a:ext cx=\"704215\" cy=\"1000125\"/></a:xfrm>' +\n              '<a:prstGeom prst=\"rect\"><a:avLst/></a:prstGeom>' +\n            '</pic:spPr>' +\n          '</pic:pic>' +\n        '</a:graphicData>' +\n      '</a:graphic>' +\n    '</wp:anchor>' +\n  '</w:drawing>';\n\nconst titleRunProps =\n  '<w:rPr><w:b/><w:bCs/><w:sz w:val=\"28\"/><w:szCs w:val=\"28\"/></w:rPr>';\n\nconst titleParaOoxml =\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n      '<pkg:xmlData>' +\n        '<w:document ' +\n            'xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" ' +\n            'xmlns:r=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships\" ' +\n            'xmlns:wp=\"http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing\" ' +\n            'xmlns:a=\"http://schemas.openxmlformats.org/drawingml/2006/main\" ' +\n            'xmlns:pic=\"http://schemas.openxmlformats.org/drawingml/2006/picture\">' +\n          '<w:body>' +\n            '<w:p>' +\n              '<w:pPr>' +\n                '<w:pStyle w:val=\"Normal\"/>' +\n                '<w:ind w:left=\"0\" w:right=\"6000\" w:hanging=\"0\"/>' +\n                '<w:jc w:val=\"center\"/>' +\n                titleRunProps +\n              '</w:pPr>' +\n              '<w:r>' + titleRunProps + drawingXml + '</w:r>' +\n              '<w:r>' + titleRunProps + '<w:t>CITTA\\' DI IMPERIA</w:t></w:r>' +\n            '</w:p>' +\n          '</w:body>' +\n        '</w:document>' +\n      '</pkg:xmlData>' +\n    '</pkg:part>' +\n  '</pkg:package>';\n\ntitlePara.insertOoxml(titleParaOoxml, Word.InsertLocation.replace);\nawait context.sync();\n\n// --- 2) Text fix: \"informarLe\" -> \"informarLa\" -----------------------------\nconst hits = context.document.body.search(\"informarLe\", { matchCase: true });\nhits.load(\"items\");\nawait context.sync();\n\nfor (const hit of hits.items) {\n  hit.insertText(\"informarLa\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Edit: 1) split the title paragraph's run so the \"CITTA' DI IMPERIA\" picture\n#          is its own run that precedes a separate run holding the text, and\n#       2) fix \"informarLe\" -> \"informarLa\" in the body paragraph.\n#\n# $d / $word / $app resolve to the open document / application.\n\n$d = $word.ActiveDocument\n\n# --- 1) Re-shape the first paragraph (picture + \"CITTA' DI IMPERIA\") -------\n$titlePara = $d.Paragraphs(1)\n$titleRange = $titlePara.Range\n\n# The picture keeps its original anchor/position/size - only the run\n# structure changes: the drawing becomes its own run, immediately followed\n# by a run that just carries the \"CITTA' DI IMPERIA\" text.\n$drawingXml = '<w:drawing><wp:anchor behindDoc=\"0\" distT=\"0\" distB=\"0\" distL=\"0\" distR=\"0\" simplePos=\"0\" locked=\"0\" layoutInCell=\"1\" allowOverlap=\"1\" relativeHeight=\"2\"><wp:simplePos x=\"0\" y=\"0\"/><wp:positionH relativeFrom=\"column\"><wp:posOffset>786765</wp:posOffset></wp:positionH><wp:positionV relativeFrom=\"paragraph\"><wp:posOffset>36195</wp:posOffset></wp:positionV><wp:extent cx=\"704215\" cy=\"1000125\"/><wp:effectExtent l=\"0\" t=\"0\" r=\"0\" b=\"0\"/><wp:wrapTopAndBottom/><wp:docPr id=\"1\" name=\"Picture\" descr=\"\"/><wp:cNvGraphicFramePr><a:graphicFrameLocks xmlns:a=\"http://schemas.openxmlformats.org/drawingml/2006/main\" noChangeAspect=\"1\"/></wp:cNvGraphicFramePr><a:graphic xmlns:a=\"http://schemas.openxmlformats.org/drawingml/2006/main\"><a:graphicData uri=\"http://schemas.openxmlformats.org/drawingml/2006/picture\"><pic:pic xmlns:pic=\"http://schemas.openxmlformats.org/drawingml/2006/picture\"><pic:nvPicPr><pic:cNvPr id=\"1\" name=\"Picture\" descr=\"\"/><pic:cNvPicPr><a:picLocks noChangeAspect=\"1\" noChangeArrowheads=\"1\"/></pic:cNvPicPr></pic:nvPicPr><pic:blipFill><a:blip r:embed=\"rId2\"/><a:stretch><a:fillRect/></a:stretch></pic:blipFill><pic:spPr bwMode=\"auto\"><a:xfrm><a:off x=\"0\" y=\"0\"/><a:ext cx=\"704215\" cy=\"1000125\"/></a:xfrm><a:prstGeom prst=\"rect\"><a:avLst/></a:prstGeom></pic:spPr></pic:pic></a:graphicData></a:graphic></wp:anchor></w:drawing>'\n\n$titleRunProps = '<w:rPr><w:b/><w:bCs/><w:sz w:val=\"28\"/><w:szCs w:val=\"28\"/></w:rPr>'\n\n$titleText = \"CITTA' DI IMPERIA\"\n\n$titleParaOoxml = '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n      '<w:document ' +\n          'xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" ' +\n          'xmlns:r=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships\" ' +\n          'xmlns:wp=\"http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing\" ' +\n          'xmlns:a=\"http://schemas.openxmlformats.org/drawingml/2006/main\" ' +\n          'xmlns:pic=\"http://schemas.openxmlformats.org/drawingml/2006/picture\">' +\n        '<w:body>' +\n          '<w:p>' +\n            '<w:pPr>' +\n              '<w:pStyle w:val=\"Normal\"/>' +\n              '<w:ind w:left=\"0\" w:right=\"6000\" w:hanging=\"0\"/>' +\n              '<w:jc w:val=\"center\"/>' +\n              $titleRunProps +\n            '</w:pPr>' +\n            '<w:r>' + $titleRunProps + $drawingXml + '</w:r>' +\n            '<w:r>' + $titleRunProps + '<w:t>' + $titleText + '</w:t></w:r>' +\n          '</w:p>' +\n        '</w:body>' +\n      '</w:document>' +\n    '</pkg:xmlData>' +\n  '</pkg:part>' +\n'</pkg:package>'\n\n$titleRange.InsertXML($titleParaOoxml)\n\n# --- 2) Text fix: \"informarLe\" -> \"informarLa\" -----------------------------\n$find = $d.Content.Find\n$find.Text = \"informarLe\"\n$find.Replacement.Text = \"informarLa\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n"}
